$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D12: clear the title text (becomes an empty cell)
$ws.Range("D12").Value = ""

# E12: update link to the new blog post URL
$ws.Range("E12").Value = "https://tensorflow.blog/2025/03/26/book-roadmap/"

# D33: update title from blog name to the actual post title
$ws.Range("D33").Value = "Tabular Data(정형 데이터)에서의 Noise"
